# Link the agribalyse "production" activities to their ecoinvent/ex-ante
# counterparts by filling in the short "reference product" name used for
# matching (egg / sea bass or sea bream / large trout / small trout).
#
# Each of the 4 blocks in the sheet has:
#   - a "reference product" row (column A label, column B left blank before)
#   - an Exchanges sub-table whose first data row is the "production" row,
#     where column B (the reference product name) was also left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Egg block
$ws.Range("B4").Value = "egg"
$ws.Range("B11").Value = "egg"

# Sea bass / sea bream block
$ws.Range("B15").Value = "sea bass or sea bream"
$ws.Range("B22").ClearFormats()
$ws.Range("B22").Value = "sea bass or sea bream"

# Large trout block
$ws.Range("B26").Value = "large trout"
$ws.Range("B33").Value = "large trout"

# Small trout block
$ws.Range("B37").Value = "small trout"
$ws.Range("B44").Value = "small trout"

# Reflect the author's final cursor position before saving
$ws.Range("B49").Select() | Out-Null
